# "moved enum to week 9"
# Week 8 (row 8) previously covered "Interfaces, Object base class, enumerated
# types". Week 9 (row 9) previously covered "Composing data structures,
# generics". The enumerated-types topic is moved from week 8 to week 9, and
# phrasing is tightened accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Object base class, Generics"
$ws.Range("C9").Value = "Enumerated types, composing data structures"

# Reflect the author's active-cell selection ending on the updated week 9 row.
$ws.Range("C9").Select()
